$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = "@"
$ws.Range('D2').Value = '34.797.54'
$ws.Range('D2').Style = "Normal"
$ws.Range('E2').Value = '  -2.30%  '
$ws.Range('D3').NumberFormat = "@"
$ws.Range('D3').Value = '1.800.02'
$ws.Range('D3').Style = "Normal"
$ws.Range('E3').Value = '  -3.43%  '
$ws.Range('E4').Value = '  +0.09%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '230.66'
$ws.Range('D5').Style = "Normal"
$ws.Range('E5').Value = '  +0.17%  '
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '0.603'
$ws.Range('D6').Style = "Normal"
$ws.Range('E6').Value = '  -1.71%  '
$ws.Range('E7').Value = '  +0.24%  '
$ws.Range('D8').NumberFormat = "@"
$ws.Range('D8').Value = '39.56'
$ws.Range('D8').Style = "Normal"
$ws.Range('E8').Value = '  -5.60%  '
$ws.Range('D9').NumberFormat = "@"
$ws.Range('D9').Value = '0.321'
$ws.Range('D9').Style = "Normal"
$ws.Range('E9').Value = '  +3.07%  '
$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '0.0675'
$ws.Range('D10').Style = "Normal"
$ws.Range('E10').Value = '  -3.34%  '
$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '0.0988'
$ws.Range('D11').Style = "Normal"
$ws.Range('E11').Value = '  -1.63%  '
$ws.Range('D12').NumberFormat = "@"
$ws.Range('D12').Value = '2.059.93'
$ws.Range('D12').Style = "Normal"
$ws.Range('E12').Value = '  -3.50%  '
$ws.Range('D13').NumberFormat = "@"
$ws.Range('D13').Value = '1.798.80'
$ws.Range('D13').Style = "Normal"
$ws.Range('E13').Value = '  -3.39%  '
$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '0.655'
$ws.Range('D14').Style = "Normal"
$ws.Range('E14').Value = '  -3.67%  '
$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '10.83'
$ws.Range('D15').Style = "Normal"
$ws.Range('E15').Value = '  -6.58%  '
$ws.Range('D16').NumberFormat = "@"
$ws.Range('D16').Value = '4.55'
$ws.Range('D16').Style = "Normal"
$ws.Range('E16').Value = '  -5.08%  '
$ws.Range('D17').NumberFormat = "@"
$ws.Range('D17').Value = '34.719.78'
$ws.Range('D17').Style = "Normal"
$ws.Range('E17').Value = '  -2.47%  '
$ws.Range('D18').NumberFormat = "@"
$ws.Range('D18').Value = '68.69'
$ws.Range('D18').Style = "Normal"
$ws.Range('E18').Value = '  -2.50%  '
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '0.0₃0775'
$ws.Range('D19').Style = "Normal"
$ws.Range('E19').Value = '  -3.97%  '
$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '236.00'
$ws.Range('D20').Style = "Normal"
$ws.Range('E20').Value = '  -4.52%  '
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '11.70'
$ws.Range('D21').Style = "Normal"
$ws.Range('E21').Value = '  -4.57%  '
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '4.57'
$ws.Range('D22').Style = "Normal"
$ws.Range('E22').Value = '  -5.03%  '
$ws.Range('E23').Value = '  +0.38%  '
$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '2.22'
$ws.Range('D24').Style = "Normal"
$ws.Range('E24').Value = '  -0.49%  '
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '170.57'
$ws.Range('D25').Style = "Normal"
$ws.Range('E25').Value = '  -0.95%  '
$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '7.69'
$ws.Range('D26').Style = "Normal"
$ws.Range('E26').Value = '  -3.31%  '
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '17.14'
$ws.Range('D27').Style = "Normal"
$ws.Range('E27').Value = '  -4.36%  '
$ws.Range('E28').Value = '  -4.08%  '
$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '1.53'
$ws.Range('D29').Style = "Normal"
$ws.Range('E29').Value = '  +7.41%  '
$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '1.00'
$ws.Range('D30').Style = "Normal"
$ws.Range('E30').Value = '  -0.05%  '
$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '3.98'
$ws.Range('D31').Style = "Normal"
$ws.Range('E31').Value = '  +0.70%  '
$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '0.0545'
$ws.Range('D32').Style = "Normal"
$ws.Range('E32').Value = '  -0.06%  '
$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '3.90'
$ws.Range('D33').Style = "Normal"
$ws.Range('E33').Value = '  -4.30%  '
$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '1.74'
$ws.Range('D34').Style = "Normal"
$ws.Range('E34').Value = '  -9.01%  '
$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '1.13'
$ws.Range('D35').Style = "Normal"
$ws.Range('E35').Value = '  +3.08%  '
$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '0.671'
$ws.Range('D36').Style = "Normal"
$ws.Range('E36').Value = '  -2.63%  '
$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '89.86'
$ws.Range('D37').Style = "Normal"
$ws.Range('E37').Value = '  +0.52%  '
$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '1.28'
$ws.Range('D38').Style = "Normal"
$ws.Range('E38').Value = '  -0.06%  '
$ws.Range('B39').Value = 'Maker'
$ws.Range('C39').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '1.301.98'
$ws.Range('D39').Style = "Normal"
$ws.Range('E39').Value = '  -4.23%  '
$ws.Range('B40').Value = 'VeChain'
$ws.Range('C40').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '0.0190'
$ws.Range('D40').Style = "Normal"
$ws.Range('E40').Value = '  -3.50%  '
$ws.Range('B41').Value = 'HuobiToken'
$ws.Range('C41').Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '2.47'
$ws.Range('D41').Style = "Normal"
$ws.Range('E41').Value = '  +0.18%  '
$ws.Range('B42').Value = 'ARBITRUM'
$ws.Range('C42').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '0.955'
$ws.Range('D42').Style = "Normal"
$ws.Range('E42').Value = '  -7.73%  '
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '14.30'
$ws.Range('D43').Style = "Normal"
$ws.Range('E43').Value = '  -6.09%  '
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '2.74'
$ws.Range('D44').Style = "Normal"
$ws.Range('E44').Value = '  -2.88%  '
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '2.19'
$ws.Range('D45').Style = "Normal"
$ws.Range('E45').Value = '  -12.43%  '
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '6.13'
$ws.Range('D46').Style = "Normal"
$ws.Range('E46').Value = '  +0.22%  '
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '0.0509'
$ws.Range('D47').Style = "Normal"
$ws.Range('E47').Value = '  -2.48%  '
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '1.976.44'
$ws.Range('D48').Style = "Normal"
$ws.Range('E48').Value = '  -2.73%  '
$ws.Range('E49').Value = '  +0.27%  '
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '0.0656'
$ws.Range('D50').Style = "Normal"
$ws.Range('E50').Value = '  +5.13%  '
$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '98.82'
$ws.Range('D51').Style = "Normal"
$ws.Range('E51').Value = '  -6.01%  '
